$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New quote rows (7-11) ---

# Row 7: systems / James Clear
$ws.Range("A7").Value = "systems"
$ws.Range("B7").Value = "James Clear"
$ws.Range("C7").Value = "You do not rise to the level of your goals. You fall to the level of your systems."

# Row 8: systems / James Clear
$ws.Range("A8").Value = "systems"
$ws.Range("B8").Value = "James Clear"
$ws.Range("C8").Value = "If you want better results, then forget about setting goals. Focus on your system instead."

# Row 9: habits / James Clear
$ws.Range("A9").Value = "habits"
$ws.Range("B9").Value = "James Clear"
$ws.Range("C9").Value = "With outcome-based habits, the focus is on what you want to achieve. With identity-based habits, the focus is on who you wish to become."

# Row 10: habits / Jim Ryun
$ws.Range("C10").Value = "Motivation is what gets you started. Habit is what keeps you going."
$ws.Range("A10").Value = "habits"
$ws.Range("B10").Value = "Jim Ryun"

# Row 11: consistency / Anthony Robbins
$ws.Range("C11").Value = "It's not what we do once in a while that shapes our lives. It's what we do consistently."
$ws.Range("B11").Value = "Anthony Robbins"
$ws.Range("A11").Value = "consistency"

# Row heights to match the new content
$ws.Rows.Item(7).RowHeight = 18
$ws.Rows.Item(8).RowHeight = 18
$ws.Rows.Item(9).RowHeight = 18
$ws.Rows.Item(10).RowHeight = 17
$ws.Rows.Item(11).RowHeight = 17

# Widen column C to fit the new, longer quotes
$ws.Columns.Item(3).ColumnWidth = 127.5

# Move the selection to the next empty row, as Excel would after data entry
$ws.Range("A12").Select()
